$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 67
$ws.Range("D2").Value = 0.01752
$ws.Range("E2").Value = 0.00876
$ws.Range("F2").Value = 1.63116

# Row 3
$ws.Range("C3").Value = 44.65
$ws.Range("D3").Value = 0.04865
$ws.Range("E3").Value = 0.02433
$ws.Range("F3").Value = 2.6725

# Row 4
$ws.Range("C4").Value = 26.3
$ws.Range("D4").Value = 0.07103
$ws.Range("E4").Value = 0.01776
$ws.Range("F4").Value = 2.58436

# Row 5
$ws.Range("C5").Value = 36.55
$ws.Range("D5").Value = 0.03596
$ws.Range("E5").Value = 0.00899
$ws.Range("F5").Value = 2.16698

# Row 6
$ws.Range("C6").Value = 20.45
$ws.Range("D6").Value = 0.04779
$ws.Range("E6").Value = 0.00797
$ws.Range("F6").Value = 1.96758

# Row 7
$ws.Range("C7").Value = 19.1
$ws.Range("D7").Value = 0.08927
$ws.Range("E7").Value = 0.01488
$ws.Range("F7").Value = 2.66363

# Row 8
$ws.Range("C8").Value = 12.25
$ws.Range("D8").Value = 0.0957
$ws.Range("E8").Value = 0.01196
$ws.Range("F8").Value = 2.32252

# Row 9
$ws.Range("C9").Value = 16.2
$ws.Range("D9").Value = 0.05457
$ws.Range("E9").Value = 0.00682
$ws.Range("F9").Value = 1.85404

# Row 10
$ws.Range("C10").Value = 13.55
$ws.Range("D10").Value = 0.04845
$ws.Range("E10").Value = 0.00485
$ws.Range("F10").Value = 1.55669

# Row 11
$ws.Range("C11").Value = 9.95
$ws.Range("D11").Value = 0.13621
$ws.Range("E11").Value = 0.01362
$ws.Range("F11").Value = 2.65155
